$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "target total amount" (F column) values for rows 17-21
$ws.Range("F17").Value = 8000000
$ws.Range("F18").Value = 8500000
$ws.Range("F19").Value = 9000000
$ws.Range("F20").Value = 9500000
$ws.Range("F21").Value = 10000000

# Extend the shared "continued principal" formula (B column) from B16 down to B21,
# mirroring the existing pattern B{row} = F{row}-A{row}
$ws.Range("B17:B21").Formula = "=F17-A17"

# Move the active selection to C4
$ws.Range("C4").Select() | Out-Null
